$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New E-column values (trailing zeroes removed) for rows 2..192
$newValues = @(
    991,9911,99111,991111,991112,9911121,99111211,99111212,9911122,991113,991114,991115,991116,99112,99113,
    991131,9911311,99113111,99113112,9911312,991132,99114,991141,991142,99115,991151,991152,991153,991154,991156,
    9911561,9911562,991157,991158,991159,99115999,99116,991161,9911611,9911612,991162,9911621,9911622,9911623,9911624,
    9911625,991163,9911631,9911632,992,9921,99211,992111,992112,9921121,99211211,99211212,9921122,992113,992114,
    992115,992116,99212,99213,992131,9921311,99213111,99213112,9921312,992132,99214,992141,992142,99215,992151,
    992152,992153,992154,992155,9921551,9921552,992156,992157,992158,992159,9922,99221,992211,992212,99222,
    992221,992222,992223,992224,992225,9923,99231,992311,992312,993,9931,99311,993111,993112,993113,
    9931131,9931132,993114,993115,993116,99312,993121,993122,9931221,9931222,993123,993124,993125,993126,99313,
    993131,9931311,9931312,993132,993133,993134,993135,99314,993141,9931411,9931412,993142,9931421,9931422,993143,
    9931431,9931432,993144,993145,993146,9931461,9931462,993147,99315,9932,994,9941,99411,994111,994112,
    994113,9941131,9941132,994114,994115,994116,99412,994121,994122,9941221,9941222,994123,994124,994125,994126,
    99413,994131,9941311,9941312,994132,994133,994134,994135,99414,994141,9941411,9941412,994142,9941421,9941422,
    994143,9941431,9941432,994144,994145,994146,9941461,9941462,994147,99415,9942
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $newValues[$i]
}

# Clear the custom number format style previously applied to E2:E192,
# restoring the default (General) style so the unused style entry drops out.
$ws.Range("E2:E192").Style = "Normal"

# Update the active selection to match the saved view state (G6).
[void]$ws.Range("G6").Select()
